$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18 - new PBL item: UI test debugging should not be interrupted
$ws.Range("A18").Value = 17
$ws.Range("B18").Value = "Als Entwickler möchte ich die UI-Tests debuggen können, ohne dass diese unterbrochen werden"
$ws.Range("C18").Value = "Akzeptanzkriterien: Die UI-Tests brechen nicht mehr ab, wenn die Maus bewegt wird,."
$ws.Range("C18").WrapText = $true
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = "Unresolved"
$ws.Range("G18").Value = 17
$ws.Rows.Item(18).RowHeight = 45

# Row 19 - new PBL item: UI test speed
$ws.Range("A19").Value = 18
$ws.Range("B19").Value = "Als Entwickler möchte ich, dass die UI-Tests in einer angemessenen Geschwindigkeit laufen"
$ws.Range("C19").Value = "Akzeptannzkriterien: Die UI-Tests laufen wesentlich schneller als zuvor"
$ws.Range("C19").WrapText = $true

# Story points for row 19 was entered as text "0.5" (not a number) - force a
# text cell, then put the style back to Normal so no stray numeric format
# lingers on the cell.
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.5"
$ws.Range("D19").Style = "Normal"

$ws.Range("E19").Value = "Unresolved"
$ws.Range("G19").Value = 18
$ws.Rows.Item(19).RowHeight = 45

# Reflect the new scroll position / selection from the edit session
$ws.Range("B10").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 10
